$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.940.53"
$ws.Range("E2").Value = "  -1.17%  "

$ws.Range("D3").Value = "1.818.22"
$ws.Range("E3").Value = "  -0.08%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.03"
$ws.Range("E5").Value = "  -1.06%  "

$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4639"
$ws.Range("E7").Value = "  -0.63%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3707"
$ws.Range("E8").Value = "  -1.76%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07363"
$ws.Range("E9").Value = "  -0.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8726"
$ws.Range("E10").Value = "  +0.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.47"
$ws.Range("E11").Value = "  -0.67%  "

$ws.Range("D12").Value = "1.832.67"
$ws.Range("E12").Value = "  +0.66%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.356"
$ws.Range("E13").Value = "  -1.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07105"
$ws.Range("E14").Value = "  +0.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.519"
$ws.Range("E15").Value = "  -2.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.72"
$ws.Range("E16").Value = "  -0.81%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008735"
$ws.Range("E18").Value = "  -0.35%  "

$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.73"
$ws.Range("E20").Value = "  -1.32%  "

$ws.Range("D21").Value = "26.956.06"
$ws.Range("E21").Value = "  -1.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.322"
$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("E23").Value = "  -3.21%  "

$ws.Range("D24").Value = "2.064.41"
$ws.Range("E24").Value = "  +0.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.902"
$ws.Range("E25").Value = "  -1.94%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.08"
$ws.Range("E26").Value = "  +0.48%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.43"
$ws.Range("E27").Value = "  -0.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.148"
$ws.Range("E28").Value = "  -4.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.311"
$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.69"
$ws.Range("E30").Value = "  -1.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08908"
$ws.Range("E31").Value = "  -0.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7597"
$ws.Range("E32").Value = "  -3.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.156"
$ws.Range("E33").Value = "  -2.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.467"
$ws.Range("E34").Value = "  -1.34%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.921"
$ws.Range("E35").Value = "  +0.18%  "

$ws.Range("E36").Value = "  +0.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.096"
$ws.Range("E37").Value = "  +0.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01958"
$ws.Range("E38").Value = "  -0.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05258"
$ws.Range("E39").Value = "  +0.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.943"
$ws.Range("E40").Value = "  +2.11%  "

$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5344"
$ws.Range("E42").Value = "  +0.57%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.379"
$ws.Range("E43").Value = "  +0.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1665"
$ws.Range("E44").Value = "  -1.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.458"
$ws.Range("E45").Value = "  -1.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4953"
$ws.Range("E46").Value = "  -2.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.30"
$ws.Range("E47").Value = "  -1.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.678"
$ws.Range("E48").Value = "  +0.53%  "

$ws.Range("E49").Value = "  -0.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.31"
$ws.Range("E50").Value = "  -2.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06284"
$ws.Range("E51").Value = "  -0.79%  "
